$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "CHANNEL CODE"
$ws.Range("F2").Value = "CHANNEL NAME"
$ws.Range("E3").Value = "DEPT"
$ws.Range("F3").Value = "DEPARTMENT STORE"
$ws.Range("E2:F2").Font.Bold = $true
$ws.Range("E:F").AutoFit()
